$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New daily rows to append after the last existing row (282), continuing the
# date series in column A and repeating the same accumulated totals found in
# the last row for columns B/C/D.
$newRows = @(
    @("09-10-2021", 17537, 737, 521),
    @("10-10-2021", 17537, 737, 521),
    @("11-10-2021", 17537, 737, 521),
    @("12-10-2021", 17537, 737, 521),
    @("13-10-2021", 17537, 737, 521),
    @("14-10-2021", 17537, 737, 521)
)

$startRow = 283
$scratchRow = 1000

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds a date-like label (dd-mm-yyyy) that must be stored as
    # plain text, matching every other cell in the column. Assigning the
    # string straight to .Value lets Excel's smart entry reinterpret it as a
    # real date serial (and stamp a number-format style on the cell), so we
    # build the text in a scratch cell via a text formula (guaranteed to stay
    # a string) and copy only the computed value across - this preserves the
    # shared-string/plain-text representation with no extra styling.
    $scratch = $ws.Cells.Item($scratchRow, 1)
    $scratch.Formula = '="' + $row[0] + '"'
    $scratch.Copy()
    $target = $ws.Cells.Item($r, 1)
    $target.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = $false
